# Fill in the "PERSON 3" evaluation sheet, mirroring the layout/formatting
# already used by "PERSON 2", then overwrite it with Person 3's own answers.

$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("PERSON 2")
$dst = $wb.Worksheets.Item("PERSON 3")

# Bring over the table structure (values + cell styles) from PERSON 2 as a
# starting point -- this also carries the borders/fills/fonts used by the
# evaluation template.
$src.Range("A1:F11").Copy($dst.Range("A1"))

# Match column widths / row heights used on the PERSON 2 sheet.
for ($i = 1; $i -le 6; $i++) {
    $dst.Columns.Item($i).ColumnWidth = $src.Columns.Item($i).ColumnWidth
}
foreach ($i in 4, 5, 8, 9, 10, 11) {
    $dst.Rows.Item($i).RowHeight = $src.Rows.Item($i).RowHeight
}

# --- Candidate info -------------------------------------------------------
$dst.Range("C2").Value = 28
$dst.Range("C3").Value = "M"
$dst.Range("C4").Value = "Masters"
$dst.Range("C5").Value = 0.39374999999999999

# --- Q1..Q4 answers ---------------------------------------------------
$dst.Range("C8").Value = "Partly"
$dst.Range("D8").Value = "Yes"
$dst.Range("E8").Value = "Yes"
$dst.Range("F8").Value = "Partly (south of region Africa or country South Africa?; future or current?)"

$dst.Range("C9").Value = "Yes"
$dst.Range("D9").Value = "No"
$dst.Range("E9").Value = "Yes"
$dst.Range("F9").Value = "Yes"

$dst.Range("C10").Value = "Partly"
$dst.Range("D10").Value = "Yes"
$dst.Range("E10").Value = "Yes"
$dst.Range("F10").Value = "Yes"

$dst.Range("C11").Value = "Only issue was uncertainty about 'currently': unclear that only 'operating' should be selected."
$dst.Range("D11").Value = "Clear question and easy to find answer. More info when hovering about time slider would be helpful."
$dst.Range("E11").Value = "Unclear to first select a subregion before being able to select country. Would be better to directly be able to select a country.  "
$dst.Range("F11").Value = "Selecting country directly would be helpful."

# --- Selections / active sheet, matching the final state of the workbook --
$src.Activate() | Out-Null
$src.Range("E11").Select() | Out-Null

$dst.Activate() | Out-Null
$dst.Range("F11").Select() | Out-Null

Write-Output "PERSON 3 sheet populated"
